$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price / volume(1h) text values.
# NumberFormat is forced to "@" (Text) on each cell right before the value
# is written so Excel keeps them as literal text (e.g. "-1.34%") instead of
# auto-converting them into numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "325.46"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.34%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.14%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.704"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.94%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08026"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.98%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.039"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "5.64%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.638"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.19%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.494"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.83%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.943"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.15%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9226"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-2.15%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1255"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.77%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1958"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.81%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.749"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "21.02%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09170"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.31%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03561"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.11%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.23%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001284"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.87%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006270"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.12%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.365"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.04%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3479"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.26%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1352"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.53%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2668"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.21%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04404"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.51%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001259"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004606"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.81%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001189"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.89%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02500"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.24%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05311"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.73%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007465"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.48%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009901"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.53%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.57%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002113"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.59%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01161"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.52%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006685"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.28%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003038"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.16%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002277"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.14%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"

# Row 7 / Row 8 coin swap (KuCoinToken <-> GateToken) - name and link columns
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
